$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.278.85'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '1.803.63'
$ws.Range('E3').Value = '  +3.04%  '
$ws.Range('D4').Value = "'" + '1.004'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'" + '339.56'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').Value = "'" + '1.000'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = "'" + '0.4637'
$ws.Range('E7').Value = '  +19.83%  '
$ws.Range('D8').Value = "'" + '0.3811'
$ws.Range('E8').Value = '  +12.57%  '
$ws.Range('D9').Value = "'" + '45.28'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = "'" + '1.157'
$ws.Range('E10').Value = '  +4.18%  '
$ws.Range('D11').Value = "'" + '0.07606'
$ws.Range('E11').Value = '  +5.70%  '
$ws.Range('D12').Value = "'" + '22.52'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = "'" + '6.359'
$ws.Range('E14').Value = '  +3.05%  '
$ws.Range('D15').Value = "'" + '7.563'
$ws.Range('E15').Value = '  +6.77%  '
$ws.Range('D16').Value = '1.807.70'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').Value = "'" + '0.00001095'
$ws.Range('E17').Value = '  +3.62%  '
$ws.Range('D18').Value = "'" + '0.06721'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('D19').Value = "'" + '81.61'
$ws.Range('E19').Value = '  +2.87%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  +4.60%  '
$ws.Range('D22').Value = "'" + '6.442'
$ws.Range('E22').Value = '  +4.27%  '
$ws.Range('D23').Value = '28.275.74'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').Value = "'" + '2.424'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = "'" + '20.69'
$ws.Range('E26').Value = '  +4.19%  '
$ws.Range('D27').Value = "'" + '153.27'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('D29').Value = '2.013.61'
$ws.Range('E29').Value = '  +3.10%  '
$ws.Range('D30').Value = "'" + '133.09'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('D31').Value = "'" + '1.254'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('D32').Value = "'" + '4.039'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').Value = "'" + '0.09578'
$ws.Range('E33').Value = '  +9.06%  '
$ws.Range('D34').Value = "'" + '5.854'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').Value = "'" + '0.2306'
$ws.Range('E35').Value = '  +9.58%  '
$ws.Range('D36').Value = "'" + '12.11'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = "'" + '5.282'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'" + '0.02355'
$ws.Range('E38').Value = '  +3.60%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = "'" + '0.06355'
$ws.Range('E39').Value = '  +4.18%  '
$ws.Range('D40').Value = "'" + '0.6636'
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range('E41').Value = '  +3.11%  '
$ws.Range('D42').Value = "'" + '8.385'
$ws.Range('E42').Value = '  +4.73%  '
$ws.Range('D43').Value = "'" + '1.490'
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('D44').Value = "'" + '14.15'
$ws.Range('E44').Value = '  +3.46%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = "'" + '0.6145'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('D47').Value = "'" + '3.876'
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('D48').Value = "'" + '131.00'
$ws.Range('E48').Value = '  +3.36%  '
$ws.Range('E49').Value = '  +2.64%  '
$ws.Range('D50').Value = "'" + '0.07205'
$ws.Range('E50').Value = '  +3.49%  '
$ws.Range('D51').Value = "'" + '1.180'
$ws.Range('E51').Value = '  +1.89%  '
